$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Sheet1 ("Sheet1"): append Week 24 abundance-estimate rows (158-164)
# ---------------------------------------------------------------------

# Row 158 - Suisun Bay
$ws1.Range("A158").Value = 24
$ws1.Range("B158").Value = "Suisun Bay"
$ws1.Range("C158").Value = 6
$ws1.Range("D158").Value = 24
$ws1.Range("E158").Value = 0
$ws1.Range("F158").Value = 0
$ws1.Range("G158").Value = 0
$ws1.Range("H158").Value = 105412
$ws1.Range("H158").NumberFormat = "#,##0"
$ws1.Range("I158").Value = "0*"
$ws1.Range("J158").Value = "NA"
$ws1.Range("K158").Value = "NA"

# Row 159 - Suisun Marsh
$ws1.Range("A159").Value = 24
$ws1.Range("B159").Value = "Suisun Marsh"
$ws1.Range("C159").Value = 5
$ws1.Range("D159").Value = 16
$ws1.Range("E159").Value = 5
$ws1.Range("F159").Value = 0
$ws1.Range("G159").Value = 2
$ws1.Range("H159").Value = 66643
$ws1.Range("H159").NumberFormat = "#,##0"
$ws1.Range("I159").Value = 3371
$ws1.Range("I159").NumberFormat = "#,##0"
$ws1.Range("J159").Value = 506
$ws1.Range("K159").Value = 11786
$ws1.Range("K159").NumberFormat = "#,##0"

# Row 160 - Lower Sacramento
$ws1.Range("A160").Value = 24
$ws1.Range("B160").Value = "Lower Sacramento"
$ws1.Range("C160").Value = 6
$ws1.Range("D160").Value = 24
$ws1.Range("E160").Value = 0
$ws1.Range("F160").Value = 0
$ws1.Range("G160").Value = 0
$ws1.Range("H160").Value = 94953
$ws1.Range("H160").NumberFormat = "#,##0"
$ws1.Range("I160").Value = "0*"
$ws1.Range("J160").Value = "NA"
$ws1.Range("K160").Value = "NA"

# Row 161 - Cache Slough LI
$ws1.Range("A161").Value = 24
$ws1.Range("B161").Value = "Cache Slough LI"
$ws1.Range("C161").Value = 3
$ws1.Range("D161").Value = 12
$ws1.Range("E161").Value = 0
$ws1.Range("F161").Value = 0
$ws1.Range("G161").Value = 0
$ws1.Range("H161").Value = 36871
$ws1.Range("H161").NumberFormat = "#,##0"
$ws1.Range("I161").Value = "0*"
$ws1.Range("J161").Value = "NA"
$ws1.Range("K161").Value = "NA"

# Row 162 - Sac DW Ship Channel
$ws1.Range("A162").Value = 24
$ws1.Range("B162").Value = "Sac DW Ship Channel"
$ws1.Range("C162").Value = 6
$ws1.Range("D162").Value = 24
$ws1.Range("E162").Value = 0
$ws1.Range("F162").Value = 0
$ws1.Range("G162").Value = 0
$ws1.Range("H162").Value = 73983
$ws1.Range("H162").NumberFormat = "#,##0"
$ws1.Range("I162").Value = "0*"
$ws1.Range("J162").Value = "NA"
$ws1.Range("K162").Value = "NA"

# Row 163 - Lower San Joaquin
$ws1.Range("A163").Value = 24
$ws1.Range("B163").Value = "Lower San Joaquin"
$ws1.Range("C163").Value = 4
$ws1.Range("D163").Value = 14
$ws1.Range("E163").Value = 0
$ws1.Range("F163").Value = 0
$ws1.Range("G163").Value = 0
$ws1.Range("H163").Value = 60361
$ws1.Range("H163").NumberFormat = "#,##0"
$ws1.Range("I163").Value = "0*"
$ws1.Range("J163").Value = "NA"
$ws1.Range("K163").Value = "NA"

# Row 164 - All Strata
$ws1.Range("A164").Value = 24
$ws1.Range("B164").Value = "All Strata"
$ws1.Range("C164").Value = 30
$ws1.Range("D164").Value = 114
$ws1.Range("E164").Value = 5
$ws1.Range("F164").Value = 0
$ws1.Range("G164").Value = 2
$ws1.Range("H164").Value = 438224
$ws1.Range("H164").NumberFormat = "#,##0"
$ws1.Range("I164").Value = 3371
$ws1.Range("I164").NumberFormat = "#,##0"
$ws1.Range("J164").Value = 506
$ws1.Range("K164").Value = 11786
$ws1.Range("K164").NumberFormat = "#,##0"

# ---------------------------------------------------------------------
# Sheet2 ("Sheet2"): append Week 24 -> date-range lookup row (25)
# ---------------------------------------------------------------------
$ws2.Range("A25").Value = 24
$ws2.Range("B25").Value = "November 10–14, 2025"

# ---------------------------------------------------------------------
# Update sheet view / selection state to mirror the saved workbook state
# (Sheet1 stays the active/selected tab, so select on Sheet2 first.)
# ---------------------------------------------------------------------
$ws2.Range("A25:B25").Select()

$ws1.Activate()
$ws1.Range("B160").Select()
$excel.ActiveWindow.ScrollRow = 147

# Sheet1 picked up an explicit (portrait) page setup
$ws1.PageSetup.Orientation = 1
